$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.539082
$ws.Range("H2").Value = 58.61724600000001
$ws.Range("I2").Value = 0.224220971665117
$ws.Range("J2").Value = 0.224220971665117
$ws.Range("M2").Value = 28.063402
$ws.Range("N2").Value = 84.190206
$ws.Range("O2").Value = 0.2422582722789452
$ws.Range("P2").Value = 0.2422582722789452
$ws.Range("Q2").Value = 548.3331128769642
$ws.Range("R2").Value = 4934.998015892677
$ws.Range("S2").Value = 0.05431938520429756
$ws.Range("T2").Value = 0.05431938520429756
$ws.Range("G3").Value = 19.539082
$ws.Range("H3").Value = 58.61724600000001
$ws.Range("I3").Value = 0.224220971665117
$ws.Range("J3").Value = 0.224220971665117
$ws.Range("O3").Value = 0.270119931190762
$ws.Range("P3").Value = 0.2701199311907619
$ws.Range("Q3").Value = 611.3958517354401
$ws.Range("R3").Value = 5502.562665618961
$ws.Range("S3").Value = 0.06056655343770719
$ws.Range("T3").Value = 0.06056655343770717
$ws.Range("G4").Value = 19.539082
$ws.Range("H4").Value = 58.61724600000001
$ws.Range("I4").Value = 0.224220971665117
$ws.Range("J4").Value = 0.224220971665117
$ws.Range("M4").Value = 25.48508733333334
$ws.Range("N4").Value = 76.455262
$ws.Range("O4").Value = 0.22000088322333
$ws.Range("P4").Value = 0.2200008832233299
$ws.Range("Q4").Value = 497.9552111831615
$ws.Range("R4").Value = 4481.596900648453
$ws.Range("S4").Value = 0.04932881180351898
$ws.Range("T4").Value = 0.04932881180351897
$ws.Range("G5").Value = 19.539082
$ws.Range("H5").Value = 58.61724600000001
$ws.Range("I5").Value = 0.224220971665117
$ws.Range("J5").Value = 0.224220971665117
$ws.Range("M5").Value = 31.00143166666667
$ws.Range("N5").Value = 93.004295
$ws.Range("O5").Value = 0.2676209133069629
$ws.Range("P5").Value = 0.2676209133069628
$ws.Range("Q5").Value = 605.7395154523967
$ws.Range("R5").Value = 5451.65563907157
$ws.Range("S5").Value = 0.06000622121959324
$ws.Range("T5").Value = 0.06000622121959322
$ws.Range("I6").Value = 0.3010605798326856
$ws.Range("J6").Value = 0.3010605798326856
$ws.Range("M6").Value = 28.063402
$ws.Range("N6").Value = 84.190206
$ws.Range("O6").Value = 0.2422582722789452
$ws.Range("P6").Value = 0.2422582722789452
$ws.Range("Q6").Value = 736.2446236775573
$ws.Range("R6").Value = 6626.201613098016
$ws.Range("S6").Value = 0.07293441592156388
$ws.Range("T6").Value = 0.07293441592156387
$ws.Range("I7").Value = 0.3010605798326856
$ws.Range("J7").Value = 0.3010605798326856
$ws.Range("O7").Value = 0.270119931190762
$ws.Range("P7").Value = 0.2701199311907619
$ws.Range("S7").Value = 0.08132246310865596
$ws.Range("T7").Value = 0.08132246310865594
$ws.Range("I8").Value = 0.3010605798326856
$ws.Range("J8").Value = 0.3010605798326856
$ws.Range("M8").Value = 25.48508733333334
$ws.Range("N8").Value = 76.455262
$ws.Range("O8").Value = 0.22000088322333
$ws.Range("P8").Value = 0.2200008832233299
$ws.Range("Q8").Value = 668.6024215139591
$ws.Range("R8").Value = 6017.421793625632
$ws.Range("S8").Value = 0.06623359346691869
$ws.Range("T8").Value = 0.06623359346691868
$ws.Range("I9").Value = 0.3010605798326856
$ws.Range("J9").Value = 0.3010605798326856
$ws.Range("M9").Value = 31.00143166666667
$ws.Range("N9").Value = 93.004295
$ws.Range("O9").Value = 0.2676209133069629
$ws.Range("P9").Value = 0.2676209133069628
$ws.Range("Q9").Value = 813.3239651732355
$ws.Range("R9").Value = 7319.915686559119
$ws.Range("S9").Value = 0.08057010733554713
$ws.Range("T9").Value = 0.08057010733554712
$ws.Range("G10").Value = 19.67155566666667
$ws.Range("H10").Value = 59.014667
$ws.Range("I10").Value = 0.2257411748281949
$ws.Range("J10").Value = 0.2257411748281949
$ws.Range("M10").Value = 28.063402
$ws.Range("N10").Value = 84.190206
$ws.Range("O10").Value = 0.2422582722789452
$ws.Range("P10").Value = 0.2422582722789452
$ws.Range("Q10").Value = 552.0507746390447
$ws.Range("R10").Value = 4968.456971751401
$ws.Range("S10").Value = 0.05468766699609782
$ws.Range("T10").Value = 0.05468766699609781
$ws.Range("G11").Value = 19.67155566666667
$ws.Range("H11").Value = 59.014667
$ws.Range("I11").Value = 0.2257411748281949
$ws.Range("J11").Value = 0.2257411748281949
$ws.Range("O11").Value = 0.270119931190762
$ws.Range("P11").Value = 0.2701199311907619
$ws.Range("Q11").Value = 615.5410746412133
$ws.Range("R11").Value = 5539.86967177092
$ws.Range("S11").Value = 0.06097719061151378
$ws.Range("T11").Value = 0.06097719061151377
$ws.Range("G12").Value = 19.67155566666667
$ws.Range("H12").Value = 59.014667
$ws.Range("I12").Value = 0.2257411748281949
$ws.Range("J12").Value = 0.2257411748281949
$ws.Range("M12").Value = 25.48508733333334
$ws.Range("N12").Value = 76.455262
$ws.Range("O12").Value = 0.22000088322333
$ws.Range("P12").Value = 0.2200008832233299
$ws.Range("Q12").Value = 501.3313141475282
$ws.Range("R12").Value = 4511.981827327754
$ws.Range("S12").Value = 0.04966325784207503
$ws.Range("T12").Value = 0.04966325784207502
$ws.Range("G13").Value = 19.67155566666667
$ws.Range("H13").Value = 59.014667
$ws.Range("I13").Value = 0.2257411748281949
$ws.Range("J13").Value = 0.2257411748281949
$ws.Range("M13").Value = 31.00143166666667
$ws.Range("N13").Value = 93.004295
$ws.Range("O13").Value = 0.2676209133069629
$ws.Range("P13").Value = 0.2676209133069628
$ws.Range("Q13").Value = 609.8463887771961
$ws.Range("R13").Value = 5488.617498994765
$ws.Range("S13").Value = 0.0604130593785083
$ws.Range("T13").Value = 0.06041305937850829
$ws.Range("G14").Value = 21.69639766666667
$ws.Range("H14").Value = 65.08919299999999
$ws.Range("I14").Value = 0.2489772736740025
$ws.Range("J14").Value = 0.2489772736740025
$ws.Range("M14").Value = 28.063402
$ws.Range("N14").Value = 84.190206
$ws.Range("O14").Value = 0.2422582722789452
$ws.Range("P14").Value = 0.2422582722789452
$ws.Range("Q14").Value = 608.8747296715286
$ws.Range("R14").Value = 5479.872567043758
$ws.Range("S14").Value = 0.06031680415698595
$ws.Range("T14").Value = 0.06031680415698594
$ws.Range("G15").Value = 21.69639766666667
$ws.Range("H15").Value = 65.08919299999999
$ws.Range("I15").Value = 0.2489772736740025
$ws.Range("J15").Value = 0.2489772736740025
$ws.Range("O15").Value = 0.270119931190762
$ws.Range("P15").Value = 0.2701199311907619
$ws.Range("Q15").Value = 678.9002436758533
$ws.Range("R15").Value = 6110.102193082679
$ws.Range("S15").Value = 0.06725372403288507
$ws.Range("T15").Value = 0.06725372403288506
$ws.Range("G16").Value = 21.69639766666667
$ws.Range("H16").Value = 65.08919299999999
$ws.Range("I16").Value = 0.2489772736740025
$ws.Range("J16").Value = 0.2489772736740025
$ws.Range("M16").Value = 25.48508733333334
$ws.Range("N16").Value = 76.455262
$ws.Range("O16").Value = 0.22000088322333
$ws.Range("P16").Value = 0.2200008832233299
$ws.Range("Q16").Value = 552.9345893537296
$ws.Range("R16").Value = 4976.411304183566
$ws.Range("S16").Value = 0.05477522011081729
$ws.Range("T16").Value = 0.05477522011081728
$ws.Range("G17").Value = 21.69639766666667
$ws.Range("H17").Value = 65.08919299999999
$ws.Range("I17").Value = 0.2489772736740025
$ws.Range("J17").Value = 0.2489772736740025
$ws.Range("M17").Value = 31.00143166666667
$ws.Range("N17").Value = 93.004295
$ws.Range("O17").Value = 0.2676209133069629
$ws.Range("P17").Value = 0.2676209133069628
$ws.Range("Q17").Value = 672.6193896759927
$ws.Range("R17").Value = 6053.574507083935
$ws.Range("S17").Value = 0.06663152537331418
$ws.Range("T17").Value = 0.06663152537331417
